$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "UWE" (row 64, col A) used to read "UG1" in the shared-string table;
# restrict it to a single option by renaming it to "UGG1".
$ws.Range("A64").Value = "UGG1"

# "CCC" (row 64, col B) goes from allowing 100 students down to just 30.
$ws.Range("B64").Value = 30

# Reflect the scrolled/selected viewport state: the sheet view now shows
# row 55 onward with B56 as the active/selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B56").Select()
